# Bump the published term version and its publication date on the
# "Metadata" sheet of the ValueSet workbook.
#
#   Property | Value
#   ...
#   Version  | 1.0.0 -> 1.1.0          (A3 / B3)
#   ...
#   Date     | 2023-06-07T11:52:14+02:00 -> 2023-07-10T23:08:03+02:00  (A8 / B8)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

if ($ws.Name -ne "Metadata") {
    $ws = $wb.Worksheets.Item("Metadata")
}

$ws.Range("B3").Value = "1.1.0"
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
